# Update the LR-pairs data with newly recomputed TPM-based values.
# The table now covers the full 3x3 combination of sending/target clusters
# (ECs, FAPs, MuSCs) for the Wnt9a -> Fzd4 edge, replacing the previous
# partial set that also referenced "Resolving-Mac". The new table has one
# more data row (9 -> 10 rows total, i.e. through row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Wnt9a", "Fzd4", "ECs", 1, 0.3333333333333333, 0.103667, 0.311001, 0.01587706838992035, 0.01587706838992035, 3, 1, 34.97976933333334, 104.939308, 0.4352965780925344, 0.4352965780925344, 3.626247747478667, 32.63622972730801, 0.006911233540273475, 0.006911233540273473),
    @("ECs", "Wnt9a", "Fzd4", "FAPs", 1, 0.3333333333333333, 0.103667, 0.311001, 0.01587706838992035, 0.01587706838992035, 3, 1, 20.343383, 61.03014900000001, 0.2531579017099818, 0.2531579017099818, 2.108937485461, 18.980437369149, 0.004019405318898116, 0.004019405318898116),
    @("ECs", "Wnt9a", "Fzd4", "MuSCs", 1, 0.3333333333333333, 0.103667, 0.311001, 0.01587706838992035, 0.01587706838992035, 3, 1, 25.035323, 75.105969, 0.3115455201974837, 0.3115455201974837, 2.595336829441, 23.358031464969, 0.004946429530748762, 0.00494642953074876),
    @("FAPs", "Wnt9a", "Fzd4", "ECs", 3, 1, 3.232908333333333, 9.698725, 0.4951344854840667, 0.4951344854840667, 3, 1, 34.97976933333334, 104.939308, 0.4352965780925344, 0.4352965780925344, 113.0863877758111, 1017.7774899823, 0.2155303472268219, 0.2155303472268219),
    @("FAPs", "Wnt9a", "Fzd4", "FAPs", 3, 1, 3.232908333333333, 9.698725, 0.4951344854840667, 0.4951344854840667, 3, 1, 20.343383, 61.03014900000001, 0.2531579017099818, 0.2531579017099818, 65.76829242889167, 591.9146318600251, 0.1253472074093978, 0.1253472074093978),
    @("FAPs", "Wnt9a", "Fzd4", "MuSCs", 3, 1, 3.232908333333333, 9.698725, 0.4951344854840667, 0.4951344854840667, 3, 1, 25.035323, 75.105969, 0.3115455201974837, 0.3115455201974837, 80.93690435439167, 728.4321391895249, 0.154256930847847, 0.154256930847847),
    @("MuSCs", "Wnt9a", "Fzd4", "ECs", 3, 1, 3.192778666666667, 9.578336, 0.4889884461260129, 0.4889884461260129, 3, 1, 34.97976933333334, 104.939308, 0.4352965780925344, 0.4352965780925344, 111.6826612923876, 1005.143951631488, 0.2128549973254391, 0.212854997325439),
    @("MuSCs", "Wnt9a", "Fzd4", "FAPs", 3, 1, 3.192778666666667, 9.578336, 0.4889884461260129, 0.4889884461260129, 3, 1, 20.343383, 61.03014900000001, 0.2531579017099818, 0.2531579017099818, 64.95191925022935, 584.5672732520641, 0.1237912889816859, 0.1237912889816859),
    @("MuSCs", "Wnt9a", "Fzd4", "MuSCs", 3, 1, 3.192778666666667, 9.578336, 0.4889884461260129, 0.4889884461260129, 3, 1, 25.035323, 75.105969, 0.3115455201974837, 0.3115455201974837, 79.93224518750934, 719.3902066875841, 0.1523421598188879, 0.1523421598188879)
)

# The new data adds one extra row (row 10), so make sure any stale
# trailing content from the old 9-row table is cleared first.
$ws.Range("A2:T10").ClearContents()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
